$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need an explicit text
# format first, otherwise Excel auto-converts the text into a numeric value
# (losing formatting such as trailing zeros, e.g. "0.0760" -> 0.076).
$textCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D13", "D14", "D16", "D19", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D31", "D32", "D34", "D35", "D37", "D38", "D41", "D43", "D44", "D48", "D50", "D51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '37.100.47'
$ws.Range('E2').Value = '  +4.76%  '
$ws.Range('D3').Value = '1.918.70'
$ws.Range('E3').Value = '  +1.70%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '249.21'
$ws.Range('E5').Value = '  +1.19%  '
$ws.Range('D6').Value = '0.687'
$ws.Range('E6').Value = '  -0.82%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '47.55'
$ws.Range('E8').Value = '  +9.79%  '
$ws.Range('D9').Value = '0.376'
$ws.Range('E9').Value = '  +5.83%  '
$ws.Range('D10').Value = '58.06'
$ws.Range('E10').Value = '  +5.89%  '
$ws.Range('D11').Value = '0.0760'
$ws.Range('E11').Value = '  +2.10%  '
$ws.Range('E12').Value = '  +1.49%  '
$ws.Range('D13').Value = '15.56'
$ws.Range('E13').Value = '  +12.40%  '
$ws.Range('D14').Value = '0.821'
$ws.Range('E14').Value = '  +6.49%  '
$ws.Range('D15').Value = '2.198.47'
$ws.Range('E15').Value = '  +1.83%  '
$ws.Range('D16').Value = '5.13'
$ws.Range('E16').Value = '  +2.18%  '
$ws.Range('D17').Value = '1.917.95'
$ws.Range('E17').Value = '  +1.89%  '
$ws.Range('D18').Value = '37.182.49'
$ws.Range('E18').Value = '  +5.10%  '
$ws.Range('D19').Value = '74.64'
$ws.Range('E19').Value = '  +1.50%  '
$ws.Range('D20').Value = '0.0₃0854'
$ws.Range('E20').Value = '  +3.34%  '
$ws.Range('D21').Value = '13.62'
$ws.Range('E21').Value = '  +6.16%  '
$ws.Range('D22').Value = '250.74'
$ws.Range('E22').Value = '  +2.40%  '
$ws.Range('D23').Value = '5.15'
$ws.Range('E23').Value = '  +0.33%  '
$ws.Range('E24').Value = '  +0.15%  '
$ws.Range('D25').Value = '2.50'
$ws.Range('E25').Value = '  -5.09%  '
$ws.Range('D26').Value = '167.56'
$ws.Range('E26').Value = '  +1.17%  '
$ws.Range('D27').Value = '2.10'
$ws.Range('E27').Value = '  -3.06%  '
$ws.Range('D28').Value = '8.80'
$ws.Range('E28').Value = '  +2.01%  '
$ws.Range('D29').Value = '18.69'
$ws.Range('E29').Value = '  +2.16%  '
$ws.Range('E30').Value = '  +0.39%  '
$ws.Range('D31').Value = '4.55'
$ws.Range('E31').Value = '  +5.89%  '
$ws.Range('D32').Value = '0.0609'
$ws.Range('E32').Value = '  +2.16%  '
$ws.Range('E33').Value = '  +26.28%  '
$ws.Range('D34').Value = '4.29'
$ws.Range('E34').Value = '  +2.71%  '
$ws.Range('D35').Value = '1.90'
$ws.Range('E35').Value = '  +0.60%  '
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('D37').Value = '19.21'
$ws.Range('E37').Value = '  +38.95%  '
$ws.Range('D38').Value = '0.890'
$ws.Range('E38').Value = '  +3.80%  '
$ws.Range('E39').Value = '  -0.27%  '
$ws.Range('E40').Value = '  -0.25%  '
$ws.Range('D41').Value = '105.47'
$ws.Range('E41').Value = '  +7.88%  '
$ws.Range('E42').Value = '  +2.66%  '
$ws.Range('D43').Value = '17.55'
$ws.Range('E43').Value = '  +1.55%  '
$ws.Range('D44').Value = '2.91'
$ws.Range('E44').Value = '  +21.06%  '
$ws.Range('E45').Value = '  +1.83%  '
$ws.Range('D46').Value = '1.348.72'
$ws.Range('E46').Value = '  +1.81%  '
$ws.Range('E47').Value = '  +0.58%  '
$ws.Range('D48').Value = '0.0836'
$ws.Range('E48').Value = '  +3.19%  '
$ws.Range('E49').Value = '  +2.29%  '
$ws.Range('D50').Value = '6.40'
$ws.Range('E50').Value = '  +1.87%  '
$ws.Range('D51').Value = '3.76'
$ws.Range('E51').Value = '  +13.00%  '
